# Remove the stray "▼" marker run that precedes the table caption
# "表 5-1-1 功能性需求表" (the paragraph right above the functional
# requirements table). The "▼" lives in its own <w:r> run, so a plain
# Find/Replace with empty replacement text removes that run entirely
# (Word collapses/drops runs that end up with no text).

$d = $word.ActiveDocument

$marker = [char]0x25BC   # ▼  (BLACK DOWN-POINTING TRIANGLE)

$found = $d.Content.Find.Execute(
    $marker,   # FindText
    $true,     # MatchCase
    $false,    # MatchWholeWord
    $false,    # MatchWildcards
    $false,    # MatchSoundsLike
    $false,    # MatchAllWordForms
    $true,     # Forward
    1,         # Wrap (wdFindContinue)
    $false,    # Format
    "",        # ReplaceWith
    2          # Replace (wdReplaceAll)
)

Write-Host "Removed marker:" $found
